$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the "Here is a test
#    instruction." paragraph to the end of the "Author: Eric Lynch"
#    paragraph.
# ------------------------------------------------------------------

# Remove the existing bookmark (currently sitting right after the text
# of the "Here is a test instruction." paragraph).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-create it at the end of the "Author: Eric Lynch" paragraph, taking
# care to preserve that paragraph's existing w14:paraId/w:rsidR markup
# (InsertXML replaces the whole target range, so the original paragraph
# mark attributes are reproduced explicitly here).
$authorPara = $d.Paragraphs(3)
$authorXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="32617B76" w14:textId="745C9578" w:rsidR="00EB35DC" w:rsidRDefault="00EB35DC"><w:r><w:t>Author: Eric Lynch</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$authorPara.Range.InsertXML($authorXml) | Out-Null

# ------------------------------------------------------------------
# 2. Add a new paragraph of baseball-themed filler text right after
#    "Here is a test instruction.", with "choke" flagged by grammar
#    proofing marks (w:proofErr gramStart/gramEnd).
# ------------------------------------------------------------------

$testPara = $d.Paragraphs(5)
$testPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs(6)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">First base hitter designated hitter helmet strike zone first baseman cracker jack. Forkball field mustard ball cup of coffee curve team. Grand slam suicide squeeze batting average third base hardball peanuts sidearm suicide squeeze. Shutout hardball shutout sweep assist tag runs pickoff. Backstop can of corn interleague double switch on-base percentage can of corn helmet run. Run batted in knuckleball grand slam off-speed foul, error first baseman </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>choke</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> up.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml) | Out-Null
